# Generate Report for Handoff
# Regenerates the localization-status report after a new handoff of the
# "2f0c4452-87e5-450a-aff9-ca34f095f647" source file: it moves from
# "In Translation" to "Ready for handoff", while
# "d4dfe3a8-3bdf-4884-a18d-bea7584d36ec" stays "In Translation".
#
# Row 2 on every sheet now carries the d4dfe3a8 file, row 3 the 2f0c4452
# file (the two rows trade places relative to the previous report).

$wb = $excel.ActiveWorkbook

function Set-LinkCell($ws, $cellRef, $newText) {
    # Update both the cell's own text and the hyperlink's display text
    # (a file-name cell whose value is also the hyperlink caption) so the
    # two stay in sync, the way the report generator writes them.
    $ws.Range($cellRef).Value = $newText
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address(0, 0) -eq $cellRef) {
            $hl.TextToDisplay = $newText
            return
        }
    }
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 -> d4dfe3a8 file, still "In Translation" as of the earlier handoff.
Set-LinkCell $ov "A2" "d4dfe3a8-3bdf-4884-a18d-bea7584d36ec.md"
$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"
$ov.Range("D2").Value = "2016-13-20 10:13:26"

# Row 3 -> 2f0c4452 file, just became "Ready for handoff".
Set-LinkCell $ov "A3" "2f0c4452-87e5-450a-aff9-ca34f095f647.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-14-20 10:14:08"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 -> d4dfe3a8 file (unchanged "In Translation" handoff from before).
Set-LinkCell $zh "A2" "d4dfe3a8-3bdf-4884-a18d-bea7584d36ec.md"
$zh.Range("C2").Value = "In Translation"
Set-LinkCell $zh "D2" "d4dfe3a8-3bdf-4884-a18d-bea7584d36ec.4cad95e5751467f93369b43d9dea06e8d38ccfcc.zh-cn.xlf"

# Row 3 -> 2f0c4452 file, new handoff just generated.
Set-LinkCell $zh "A3" "2f0c4452-87e5-450a-aff9-ca34f095f647.md"
$zh.Range("C3").Value = "Ready for handoff"
Set-LinkCell $zh "D3" "2f0c4452-87e5-450a-aff9-ca34f095f647.bbdcdfdd908574d75b362231415b6afd8d4a4cc0.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-20 10:14:05"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2 -> d4dfe3a8 file (unchanged "In Translation" handoff from before).
Set-LinkCell $de "A2" "d4dfe3a8-3bdf-4884-a18d-bea7584d36ec.md"
$de.Range("C2").Value = "In Translation"
Set-LinkCell $de "D2" "d4dfe3a8-3bdf-4884-a18d-bea7584d36ec.4cad95e5751467f93369b43d9dea06e8d38ccfcc.de-de.xlf"

# Row 3 -> 2f0c4452 file, new handoff just generated.
Set-LinkCell $de "A3" "2f0c4452-87e5-450a-aff9-ca34f095f647.md"
$de.Range("C3").Value = "Ready for handoff"
Set-LinkCell $de "D3" "2f0c4452-87e5-450a-aff9-ca34f095f647.bbdcdfdd908574d75b362231415b6afd8d4a4cc0.de-de.xlf"
$de.Range("E3").Value = "2016-03-20 10:14:08"
